$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73, shifting existing rows 73-115 down to 74-116
$ws.Rows("73:73").Insert()

# Populate the newly inserted row 73 with the new data record
$ws.Cells.Item(73, 1).Value = 6
$ws.Cells.Item(73, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(73, 3).Value = "Metropolitana"
$ws.Cells.Item(73, 4).Value = 45176
$ws.Cells.Item(73, 5).Value = 13
$ws.Cells.Item(73, 6).Value = 100112035
$ws.Cells.Item(73, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 290
$ws.Cells.Item(73, 11).Value = 18000
$ws.Cells.Item(73, 12).Value = 20000
$ws.Cells.Item(73, 13).Value = 19103
$ws.Cells.Item(73, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(73, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(73, 16).Value = 1274
$ws.Cells.Item(73, 17).Value = 15
$ws.Cells.Item(73, 18).Value = "Hortaliza"

# Ensure D73 keeps the date number format style used throughout column D
$ws.Cells.Item(73, 4).NumberFormat = $ws.Cells.Item(74, 4).NumberFormat
